# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.725.13"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "2.303.03"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.78"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.74"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.71"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0903"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("E12").Value = "  -4.25%  "
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.987"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "2.780.69"
$ws.Range("E15").Value = "  +4.50%  "
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "2.295.12"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").Value = "42.670.62"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.31"
$ws.Range("E19").Value = "  -3.25%  "
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("E21").Value = "  -1.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.40"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.06"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("E25").Value = "  -2.72%  "
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.87"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.18"
$ws.Range("E28").Value = "  +14.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  -2.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.33"
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.24"
$ws.Range("E31").Value = "  -5.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.82"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0858"
$ws.Range("E33").Value = "  -3.97%  "
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("E36").Value = "  -3.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.55"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("E38").Value = "  -2.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.81"
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.63"
$ws.Range("E40").Value = "  -2.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.07"
$ws.Range("E41").Value = "  +10.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.58"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.13"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").Value = "1.730.91"
$ws.Range("E47").Value = "  +7.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.70"
$ws.Range("E48").Value = "  -4.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "77.31"
$ws.Range("E49").Value = "  -6.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.65"
$ws.Range("E50").Value = "  -2.94%  "
$ws.Range("E51").Value = "  -3.38%  "
